$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "[49.94527471033013, 50.16372648286298]"
$ws.Range("U2").Value = "[49.90145911349605, 50.052603747481676]"

$ws.Range("M3").Value = "[49.86561569501549, 50.16812564210174]"
$ws.Range("U3").Value = "[49.908084765869816, 50.07482461873764]"

$ws.Range("M4").Value = "[49.80139007182887, 50.117430149509794]"
$ws.Range("U4").Value = "[49.89960319885188, 50.06628971669434]"

$ws.Range("M5").Value = "[49.87809810444453, 50.1715834110755]"
$ws.Range("U5").Value = "[49.89504896268212, 50.05127630215019]"

$ws.Range("M6").Value = "[49.91152778027285, 50.188747829221185]"
$ws.Range("U6").Value = "[49.899205656629924, 50.078338160269745]"

$ws.Range("M7").Value = "[49.77665002128777, 50.02577804880564]"
$ws.Range("U7").Value = "[49.91390865913977, 50.06776402502242]"

$ws.Range("M8").Value = "[49.87667023090828, 50.16353382049697]"
$ws.Range("U8").Value = "[49.94984832387235, 50.104323299967184]"

$ws.Range("M9").Value = "[49.91118246994455, 50.12591264381588]"
$ws.Range("U9").Value = "[49.94373607699038, 50.09084498749272]"

$ws.Range("M10").Value = "[49.81341277432795, 50.12100662682321]"
$ws.Range("U10").Value = "[49.82471374825896, 49.9924911948076]"

$ws.Range("M11").Value = "[49.96075630577617, 50.29203820729022]"
$ws.Range("U11").Value = "[49.847570840443524, 50.019651041647634]"

$ws.Range("M12").Value = "[49.90061077025134, 50.15329473479398]"
$ws.Range("U12").Value = "[49.92545543012937, 50.074543993927016]"

$ws.Range("M13").Value = "[49.99300163754307, 50.25979238005417]"
$ws.Range("U13").Value = "[49.98040141297026, 50.152375550761654]"
